$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F column (time_taken) timestamps on the "data" sheet ---
$timestamps = @(
    "2021-10-05 14:20:04.247149",
    "2021-10-05 14:20:04.247157",
    "2021-10-05 14:20:04.247161",
    "2021-10-05 14:20:04.247164",
    "2021-10-05 14:20:04.247167",
    "2021-10-05 14:20:04.247170",
    "2021-10-05 14:20:04.247173",
    "2021-10-05 14:20:04.247175",
    "2021-10-05 14:20:04.247179",
    "2021-10-05 14:20:04.247181",
    "2021-10-05 14:20:04.247185",
    "2021-10-05 14:20:04.247187",
    "2021-10-05 14:20:04.247190",
    "2021-10-05 14:20:04.247193",
    "2021-10-05 14:20:04.247196",
    "2021-10-05 14:20:04.247198",
    "2021-10-05 14:20:04.247202",
    "2021-10-05 14:20:04.247205",
    "2021-10-05 14:20:04.247208",
    "2021-10-05 14:20:04.247211",
    "2021-10-05 14:20:04.247214",
    "2021-10-05 14:20:04.247216",
    "2021-10-05 14:20:04.247219",
    "2021-10-05 14:20:04.247222",
    "2021-10-05 14:20:04.247225",
    "2021-10-05 14:20:04.247228",
    "2021-10-05 14:20:04.247231",
    "2021-10-05 14:20:04.247234",
    "2021-10-05 14:20:04.247236",
    "2021-10-05 14:20:04.247239",
    "2021-10-05 14:20:04.247242"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add the new "metadata" sheet right after "data" ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$newSheet.Name = "metadata"

# Copy the header formatting (bold, border, centered) from the data sheet
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)

# Copy the style of the first index column cell too
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Early onset dementia (encompassing fronto-temporal dementia and prion disease)"
$newSheet.Range("C2").Value = 265
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.48"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2019-06-20T15:15:01.659131Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:04.243698"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/265/?format=json"
